$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 74
$srcRow = 73

# Carry the row-above formatting onto the new row (column A has the custom
# date number format/border/alignment applied via style index 2) by copying
# each source cell directly onto its destination counterpart.
for ($col = 1; $col -le 10; $col++) {
    $ws.Cells.Item($srcRow, $col).Copy($ws.Cells.Item($newRow, $col))
}

$ws.Cells.Item($newRow, 1).Value = 45630
$ws.Cells.Item($newRow, 2).Value = 116.4121952
$ws.Cells.Item($newRow, 3).Value = 0.00170247
$ws.Cells.Item($newRow, 4).Value = 0.008850780000000001
$ws.Cells.Item($newRow, 5).Value = 0.06933635
$ws.Cells.Item($newRow, 6).Value = 12792.90181321
$ws.Cells.Item($newRow, 7).Value = 465.80531254
$ws.Cells.Item($newRow, 8).Value = 0.24
$ws.Cells.Item($newRow, 9).Value = 1.7904431
$ws.Cells.Item($newRow, 10).Value = 485.38834923

$wb.Save()
